# Update "想去人数" (column F) figures across sheets, as published at
# gh-pages commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5129
$ws1.Range("F6").Value = 5129
$ws1.Range("F7").Value = 112
$ws1.Range("F8").Value = 132
$ws1.Range("F11").Value = 1159
$ws1.Range("F13").Value = 4987
$ws1.Range("F14").Value = 22
$ws1.Range("F16").Value = 78
$ws1.Range("F17").Value = 212
$ws1.Range("F18").Value = 224
$ws1.Range("F19").Value = 99
$ws1.Range("F20").Value = 244
$ws1.Range("F21").Value = 3772
$ws1.Range("F24").Value = 3668
$ws1.Range("F25").Value = 176
$ws1.Range("F26").Value = 167
$ws1.Range("F28").Value = 211
$ws1.Range("F30").Value = 203
$ws1.Range("F32").Value = 108
$ws1.Range("F35").Value = 137
$ws1.Range("F36").Value = 6512
$ws1.Range("F37").Value = 1037
$ws1.Range("F38").Value = 490
$ws1.Range("F40").Value = 971
$ws1.Range("F42").Value = 1329
$ws1.Range("F43").Value = 155
$ws1.Range("F44").Value = 653
$ws1.Range("F46").Value = 2229
$ws1.Range("F49").Value = 767
$ws1.Range("F50").Value = 909

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 15
$ws2.Range("F7").Value = 133
$ws2.Range("F9").Value = 82

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 5129
$ws4.Range("F9").Value = 5129
$ws4.Range("F10").Value = 112
$ws4.Range("F12").Value = 132
$ws4.Range("F13").Value = 82
$ws4.Range("F15").Value = 1159
$ws4.Range("F17").Value = 4987
$ws4.Range("F18").Value = 22
$ws4.Range("F20").Value = 78
$ws4.Range("F21").Value = 212
$ws4.Range("F22").Value = 224
$ws4.Range("F23").Value = 99
$ws4.Range("F24").Value = 244
$ws4.Range("F25").Value = 3772
$ws4.Range("F26").Value = 3668
$ws4.Range("F27").Value = 176
$ws4.Range("F28").Value = 167
$ws4.Range("F29").Value = 211
$ws4.Range("F31").Value = 203
$ws4.Range("F33").Value = 108
$ws4.Range("F37").Value = 6512
$ws4.Range("F38").Value = 1037
$ws4.Range("F39").Value = 490
$ws4.Range("F42").Value = 971
$ws4.Range("F43").Value = 1329
$ws4.Range("F44").Value = 155
$ws4.Range("F45").Value = 653
$ws4.Range("F46").Value = 2229
$ws4.Range("F48").Value = 767
$ws4.Range("F49").Value = 909
